$wb = $excel.ActiveWorkbook
$wsColors = $wb.Worksheets.Item("Colors")
$wsViews = $wb.Worksheets.Item("Views")

# --- New "Y" / "?" / "N" column (K) on the Views sheet ------------------
# Shared-string creation order matters (must match Y, ?, N insert order),
# so set a "Y" cell first, then a "?" cell, then a "N" cell, before filling
# in the rest.
$wsViews.Range("K4").Value = "Y"
$wsViews.Range("K13").Value = "?"
$wsViews.Range("K11").Value = "N"

$yRows = @(5,6,7,8,9,10,15,17,18,19,22,25,29)
foreach ($r in $yRows) {
    $wsViews.Range("K$r").Value = "Y"
}

$xRows = @(12,14,16,20,21,24,28)
foreach ($r in $xRows) {
    $wsViews.Range("K$r").Value = "X"
}

$nRows = @(23,26,27)
foreach ($r in $nRows) {
    $wsViews.Range("K$r").Value = "N"
}

# Row 30 mirrors J30's "Waiting for spec" formula/style in the new K column.
$wsViews.Range("J30").Copy()
$wsViews.Range("K30").PasteSpecial(-4122)
$wsViews.Range("K30").Formula = '=IF(C30="","Waiting for spec","")'

# Extend the conditional-formatting rule that used to stop at J31 so it also
# covers the new K30 cell.
$wsViews.Range("J19:J31,K30").FormatConditions.Delete()
$cond = $wsViews.Range("J19:J31,K30").FormatConditions.Add(2, 0, 'AND($B19="Yes", $C19="")')
$cond.Interior.ColorIndex = 40

# --- View / selection changes -------------------------------------------
$wsViews.Range("L12").Select()
$wsColors.Range("A1").Select()

# Views becomes the active (front-most) tab, Colors no longer tabSelected.
$wsViews.Activate()

$excel.ActiveWindow.WindowState = -4140
